# RF001 - Autenticar Usuario (GT) : 1.0 -> 1.1
#
# The test-suite workbook has three test case blocks (TC1 rows 6-13, TC2 rows
# 16-23, TC3 rows 26-33), each a 3-step login scenario assembled from a shared
# pool of step/result sentences. This edit reshuffles which alternate-flow
# sentence belongs to which test case:
#   - TC1's Step 2 action/Step 3 action are swapped, and its Step 2 expected
#     result switches from "TJSeg system down" to "CAS system down".
#   - TC2's Step 2 action/Step 3 action are swapped, and its Step 2 expected
#     result switches from "username/password incorrect" to "TJSeg system down".
#   - TC3's Step 2 expected result switches from "CAS system down" to
#     "username/password incorrect" (its step actions stay as-is).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# TC1 block (rows 6-13): Step 2 = row 11, Step 3 = row 12
$ws.Range("B11").Value = "Usuario do Sistema preenche os campos e clica no botao entrar"
$ws.Range("D11").Value = "SYSTEM alerta que o CAS (sistema de autorizacao login-senha) esta fora do ar"
$ws.Range("B12").Value = "Usuario do Sistema seleciona um nome de usuario sugerido, digita a senha e clica no botao entrar"

# TC2 block (rows 16-23): Step 2 = row 21, Step 3 = row 22
$ws.Range("B21").Value = "Usuario do Sistema seleciona um nome de usuario sugerido, digita a senha e clica no botao entrar"
$ws.Range("D21").Value = "SYSTEM alerta que o TJSeg (sistema que fornece as permissoes de acesso e escrita) esta fora do ar"
$ws.Range("B22").Value = "Usuario do Sistema preenche os campos e clica no botao entrar"

# TC3 block (rows 26-33): Step 2 = row 31
$ws.Range("D31").Value = "SYSTEM alerta que o nome de usuario e/ou senha estao incorretos"
